$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that sits after
#    "download the "release" folder."
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. "Open a command widow." -> "Open the release folder"
#    with a new "_GoBack" bookmark re-inserted in the middle of the
#    word "release" (between "...relea" and "se folder...").
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Open a command widow.") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = "Open the release folder"

    $pStart = $target.Range.Start
    $bmPos  = $pStart + 14   # right after "Open the rele" + "a"
    $rBm    = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $rBm)
}

# ------------------------------------------------------------------
# 3. Merge the "CD into the ... directory" and
#    "Type "..." + press [Enter]" bullets into a single bullet:
#    "Click on runMe.bat"
# ------------------------------------------------------------------
$cdPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("CD into the")) {
        $cdPara = $para
        break
    }
}

if ($cdPara -ne $null) {
    $r = $cdPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = "Click on runMe.bat"
}

# Find and delete the (now orphaned) "Type “..."" bullet paragraph.
$typePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Type")) {
        $typePara = $para
        break
    }
}

if ($typePara -ne $null) {
    $typePara.Range.Delete()
}
